$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 998.1
$ws.Range("I111").Value = 1043
$ws.Range("J111").Value = 978.8570999999999
$ws.Range("K111").Value = 3129
$ws.Range("L111").Value = 2936.5713
$ws.Range("M111").Value = -62
$ws.Range("N111").Value = -9070.5713

$ws.Range("H113").Value = 2220
$ws.Range("I113").Value = 2025
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 2025
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 1229
$ws.Range("N113").Value = -9508

$ws.Range("H129").Value = 1097.2025
$ws.Range("I129").Value = 345.82352
$ws.Range("J129").Value = 1303.2258
$ws.Range("K129").Value = 1037.47056
$ws.Range("L129").Value = 3909.6774
$ws.Range("M129").Value = 3962.52944
$ws.Range("N129").Value = -13909.6774

$ws.Range("H135").Value = 18549496
$ws.Range("I135").Value = 1456.44
$ws.Range("J135").Value = 250400000
$ws.Range("K135").Value = 13107.96
$ws.Range("L135").Value = 2253600000
$ws.Range("M135").Value = -10572.96
$ws.Range("N135").Value = -2253605070

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 371.8
$ws.Range("I5").Value = 301.625
$ws.Range("J5").Value = 652.5
$ws.Range("K5").Value = 301.625
$ws.Range("L5").Value = 652.5
$ws.Range("M5").Value = -188.625
$ws.Range("N5").Value = -878.5

$ws.Range("H7").Value = 2233656.8
$ws.Range("I7").Value = 2871701.5
$ws.Range("K7").Value = 2871701.5
$ws.Range("M7").Value = -2871588.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1262757
$ws.Range("I2").Value = 168.42857
$ws.Range("K2").Value = 1010.57142
$ws.Range("M2").Value = -897.57142

$ws.Range("H4").Value = 5263342
$ws.Range("I4").Value = 5263342
$ws.Range("K4").Value = 15790026
$ws.Range("M4").Value = -15789914

$ws.Range("H6").Value = 213.41667
$ws.Range("I6").Value = 196.1
$ws.Range("J6").Value = 300
$ws.Range("K6").Value = 588.3
$ws.Range("L6").Value = 900
$ws.Range("M6").Value = -475.3
$ws.Range("N6").Value = -1126

$ws.Range("H7").Value = 1307.1428
$ws.Range("I7").Value = 250
$ws.Range("J7").Value = 2716.6667
$ws.Range("K7").Value = 750
$ws.Range("L7").Value = 8150.000100000001
$ws.Range("M7").Value = -638
$ws.Range("N7").Value = -8374.000100000001

$ws.Range("H9").Value = 995.55554
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 995.55554
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 2986.66662
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -3434.66662

$ws.Range("H10").Value = 180
$ws.Range("I10").Value = 95
$ws.Range("J10").Value = 243.75
$ws.Range("K10").Value = 285
$ws.Range("L10").Value = 731.25
$ws.Range("M10").Value = -146
$ws.Range("N10").Value = -1009.25

$ws.Range("H11").Value = 132.11111
$ws.Range("I11").Value = 142.375
$ws.Range("J11").Value = 50
$ws.Range("K11").Value = 427.125
$ws.Range("L11").Value = 150
$ws.Range("M11").Value = -287.125
$ws.Range("N11").Value = -430

$ws.Range("H12").Value = 67.28570999999999
$ws.Range("J12").Value = 61.666668
$ws.Range("L12").Value = 185.000004
$ws.Range("N12").Value = -531.000004

$ws.Range("H13").Value = 350.14285
$ws.Range("I13").Value = 287.75
$ws.Range("J13").Value = 433.33334
$ws.Range("K13").Value = 863.25
$ws.Range("L13").Value = 1300.00002
$ws.Range("M13").Value = -695.25
$ws.Range("N13").Value = -1636.00002

$ws.Range("H15").Value = 10000
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 10000
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 30000
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -30280

$ws.Range("H95").Value = 13985
$ws.Range("J95").Value = 13985
$ws.Range("L95").Value = 41955
$ws.Range("N95").Value = -46073

$ws.Range("H97").Value = 1681.6
$ws.Range("I97").Value = 1250
$ws.Range("J97").Value = 1969.3334
$ws.Range("K97").Value = 3750
$ws.Range("L97").Value = 5908.0002
$ws.Range("M97").Value = -3254
$ws.Range("N97").Value = -6900.0002

$ws.Range("H106").Value = 12533.333
$ws.Range("J106").Value = 12533.333
$ws.Range("L106").Value = 37599.999
$ws.Range("N106").Value = -39491.999

$ws.Range("H109").Value = 4756.25
$ws.Range("I109").Value = 2250
$ws.Range("J109").Value = 5114.2856
$ws.Range("K109").Value = 6750
$ws.Range("L109").Value = 15342.8568
$ws.Range("M109").Value = -5710
$ws.Range("N109").Value = -17422.8568

$ws.Range("H131").Value = 1123.9153
$ws.Range("I131").Value = 406.42856
$ws.Range("J131").Value = 1347.1333
$ws.Range("K131").Value = 1219.28568
$ws.Range("L131").Value = 4041.3999
$ws.Range("M131").Value = 3820.71432
$ws.Range("N131").Value = -14121.3999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 5000
$ws.Range("J4").Value = 5000
$ws.Range("L4").Value = 5000
$ws.Range("N4").Value = -5224

$ws.Range("H5").Value = 6334.6665
$ws.Range("I5").Value = 1004
$ws.Range("J5").Value = 9000
$ws.Range("K5").Value = 1004
$ws.Range("L5").Value = 9000
$ws.Range("M5").Value = -892
$ws.Range("N5").Value = -9224

$ws.Range("H11").Value = 4289254.5
$ws.Range("I11").Value = 15000100
$ws.Range("J11").Value = 4916
$ws.Range("K11").Value = 15000100
$ws.Range("L11").Value = 4916
$ws.Range("M11").Value = -14999961
$ws.Range("N11").Value = -5194

$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()

$ws.Range("H36").Value = 2900
$ws.Range("I36").Value = 866.6667
$ws.Range("J36").Value = 9000
$ws.Range("K36").Value = 866.6667
$ws.Range("L36").Value = 9000
$ws.Range("M36").Value = -381.6667
$ws.Range("N36").Value = -9970

$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()

$ws.Range("H80").Value = 8585561
$ws.Range("I80").Value = 13891402
$ws.Range("J80").Value = 1883445.6
$ws.Range("K80").Value = 13891402
$ws.Range("L80").Value = 1883445.6
$ws.Range("M80").Value = -13890404
$ws.Range("N80").Value = -1885441.6

$ws.Range("H83").Value = 8585561
$ws.Range("I83").Value = 13891402
$ws.Range("J83").Value = 1883445.6
$ws.Range("K83").Value = 69457010
$ws.Range("L83").Value = 9417228
$ws.Range("M83").Value = -69452018
$ws.Range("N83").Value = -9427212

$ws.Range("H118").Value = 14002.895
$ws.Range("J118").Value = 14002.895
$ws.Range("L118").Value = 14002.895
$ws.Range("N118").Value = -17316.895

$ws.Range("H119").Value = 38406.668
$ws.Range("J119").Value = 38406.668
$ws.Range("L119").Value = 38406.668
$ws.Range("N119").Value = -48082.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5386.1816
$ws.Range("I122").Value = 5040.154
$ws.Range("J122").Value = 6671.4287
$ws.Range("K122").Value = 15120.462
$ws.Range("L122").Value = 20014.2861
$ws.Range("M122").Value = -12670.462
$ws.Range("N122").Value = -24914.2861

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1769.6666
$ws.Range("I122").Value = 1788.0769
$ws.Range("J122").Value = 1650
$ws.Range("K122").Value = 5364.2307
$ws.Range("L122").Value = 4950
$ws.Range("M122").Value = -2914.2307
$ws.Range("N122").Value = -9850

Write-Output "applied changes"